$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DevLog entry for row 10 - "Added acceleration animation for smoother
# movement, working on wings movement" (Content) and the follow-up To Do note.
# Order matters for shared-string table layout: the To Do note is entered
# first so it lands at the lower shared-string index, matching the Content
# entry that follows it.
$ws.Range("I10").Value = 44986
$ws.Range("I10").NumberFormat = $ws.Range("I9").NumberFormat

$ws.Range("K10").Value = "I should do the todo tasks."
$ws.Range("J10").Value = "Added acceleration animation for smoother movement, working on wings movement"

# Reflect the scrolled view / current selection like the author left it.
$ws.Range("J10").Select()
